$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns I1 ("I0") and J1 ("IF"), copying the style used by the
# existing header cells (e.g. H1) so the new headers look consistent.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the data values for the new columns I (I0) and J (IF), rows 2-14
$data = @(
    @(8, 9),
    @(8, 8),
    @(5, 6),
    @(7, 7),
    @(3, 3),
    @(7, 7),
    @(8, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(5, 5),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
